# Update 'F' column (想去人数 / interested count) values per the diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 918
$ws.Range("F4").Value = 808
$ws.Range("F5").Value = 885
$ws.Range("F6").Value = 465
$ws.Range("F7").Value = 714
$ws.Range("F9").Value = 1318
$ws.Range("F10").Value = 739
$ws.Range("F12").Value = 563
$ws.Range("F13").Value = 190
$ws.Range("F14").Value = 57
$ws.Range("F15").Value = 1181
$ws.Range("F16").Value = 144
$ws.Range("F18").Value = 425
$ws.Range("F21").Value = 603
$ws.Range("F22").Value = 159
$ws.Range("F25").Value = 1089
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 349
$ws.Range("F3").Value = 114
$ws.Range("F5").Value = 643
$ws.Range("F7").Value = 257
$ws.Range("F11").Value = 116
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 382
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 382
$ws.Range("F3").Value = 349
$ws.Range("F4").Value = 918
$ws.Range("F6").Value = 808
$ws.Range("F7").Value = 885
$ws.Range("F8").Value = 465
$ws.Range("F9").Value = 465
$ws.Range("F10").Value = 714
$ws.Range("F12").Value = 1318
$ws.Range("F13").Value = 739
$ws.Range("F14").Value = 114
$ws.Range("F17").Value = 563
$ws.Range("F18").Value = 643
$ws.Range("F19").Value = 190
$ws.Range("F20").Value = 57
$ws.Range("F21").Value = 1181
$ws.Range("F23").Value = 144
$ws.Range("F25").Value = 425
$ws.Range("F28").Value = 257
$ws.Range("F30").Value = 603
$ws.Range("F33").Value = 116
$ws.Range("F34").Value = 116
$ws.Range("F35").Value = 159
$ws.Range("F38").Value = 1089

